# GroupAssigning.xlsx — "Add files via upload"
#
# - Bank Class (sheet2): RANDBETWEEN results are frozen to static numbers
#   (formulas removed, last computed values kept).
# - Sheet3: filled in with a new "Account" class table (mirrors the
#   existing "Bank" class table on the Bank Class sheet).

$wb = $excel.ActiveWorkbook

$wsBank      = $wb.Worksheets.Item("Bank Class")
$wsAccount   = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------
# 1. Bank Class sheet — the RANDBETWEEN() formulas (one normal + one
#    shared-formula group) are replaced by their last computed, static
#    values. Formatting on this sheet is untouched.
# ---------------------------------------------------------------------
$wsBank.Range("B3").Value = 1
$wsBank.Range("B4").Value = 3
$wsBank.Range("B5").Value = 1
$wsBank.Range("B6").Value = 2
$wsBank.Range("B7").Value = 2

# ---------------------------------------------------------------------
# 2. Sheet3 — build the new "Account" class table, mirroring the layout
#    used for "Bank" on the Bank Class sheet.
# ---------------------------------------------------------------------
$wsAccount.Range("A1").Value = "Account"

$wsAccount.Range("A2").Value = "Methods"
$wsAccount.Range("B2").Value = "Groups"

$wsAccount.Range("A3").Value = "getOwnerFullName"
$wsAccount.Range("B3").Value = 3

$wsAccount.Range("A4").Value = "addNewTransaction"
$wsAccount.Range("B4").Value = 1

$wsAccount.Range("A5").Value = "getTransactionInfo"
$wsAccount.Range("B5").Value = 2

$wsAccount.Range("A6").Value = "toString"
$wsAccount.Range("B6").Value = 3

# B7 stays blank but picks up the same centered formatting as B2:B6.
$wsAccount.Range("B2:B7").HorizontalAlignment = -4108

# Title row: centered + merged across A1:B1 (same treatment as the
# "Bank" header on the Bank Class sheet).
$wsAccount.Range("A1:B1").HorizontalAlignment = -4108
$wsAccount.Range("A1:B1").Merge() | Out-Null

# Column A width — match the "Bank" column's best-fit width.
$wsAccount.Columns.Item(1).ColumnWidth = $wsBank.Columns.Item(1).ColumnWidth

# ---------------------------------------------------------------------
# 3. View state: Sheet3 zoomed to 160%, cursor resting on D7; Bank Class
#    stays the active tab with its cursor on A10.
# ---------------------------------------------------------------------
$wsAccount.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 160
$wsAccount.Range("D7").Select() | Out-Null

$wsBank.Activate() | Out-Null
$wsBank.Range("A10").Select() | Out-Null
